{"js": "// Each long single run of text is split into multiple segments joined by\n// manual line breaks (Word's \"Text Wrapping Break\", OOXML <w:br/>), matching\n// the target diff's <w:t>...</w:t><w:br/><w:t>...</w:t> pattern. A literal\n// vertical-tab (\\v) in the replacement text is what Word turns into a\n// manual line break when it is written back into the run.\nconst replacements = [\n  [\"1.The importance of oxygen transfer in biotechnological processes.2.Oxygen transfer systems.3.Oxygen transfer and microbial respiration: joint analysis of transfer and oxygen consumption, determination of kLa and QO2 in the fermentative process.4.Oxygen transfer in agitated and aerated systems: agitation of Newtonian liquids, agitation of Newtonian liquids subjected to aeration, agitation of non-Newtonian liquids, oxygen transfer.5.Scale-up variation: criteria for scaling up, comparisons between criteria for scaling up, scale-down.6.Downstream of biotechnological products: clarification, cell disruption, bioproduct separation techniques.\", \"1.The importance of oxygen transfer in biotechnological processes.\\v2.Oxygen transfer systems.\\v3.Oxygen transfer and microbial respiration: joint analysis of transfer and oxygen consumption, determination of kLa and QO2 in the fermentative process.\\v4.Oxygen transfer in agitated and aerated systems: agitation of Newtonian liquids, agitation of Newtonian liquids subjected to aeration, agitation of non-Newtonian liquids, oxygen transfer.\\v5.Scale-up variation: criteria for scaling up, comparisons between criteria for scaling up, scale-down.\\v6.Downstream of biotechnological products: clarification, cell disruption, bioproduct separation techniques.\"],\n  [\"A nota final (NF) ser\u00e1 composta pelas m\u00e9dias M1  e M2,calculadas conforme segue:M1=P1+a1\u00d7T1M2=P2+a2\u00d7T2Em que:-P1 e P2 s\u00e3o as notas da primeira e da segunda prova escrita, respectivamente (notas de zero a dez).-T1 e T2 s\u00e3o as notas m\u00e9dias dos trabalhos (notas de zero a dez) realizados antes da primeira e da segunda prova escrita, respectivamente.-a1 e a2 s\u00e3o os fatores multiplicadores das notas m\u00e9dias dos trabalhos, a serem definidos pelo docente antes do in\u00edcio de cada turma com base nas atividades espec\u00edficas a serem propostas. Os valores ser\u00e3o \u22650,1, sendo informados aos alunos no in\u00edcio do semestre. Em todos os casos, os valores m\u00e1ximos para M1 e M2 ser\u00e3o \u201cdez\u201d, sendo desconsideradas pontua\u00e7\u00f5es superiores.O c\u00e1lculo de NF ser\u00e1 feito conforme segue:NF=(M1+2\u00d7M2)/3Ser\u00e3o aprovados os alunos que obtiverem NF maior ou igual 5,0.\", \"A nota final (NF) ser\u00e1 composta pelas m\u00e9dias M1  e M2,calculadas conforme segue:\\vM1=P1+a1\u00d7T1\\vM2=P2+a2\u00d7T2\\vEm que:\\v-P1 e P2 s\u00e3o as notas da primeira e da segunda prova escrita, respectivamente (notas de zero a dez).\\v-T1 e T2 s\u00e3o as notas m\u00e9dias dos trabalhos (notas de zero a dez) realizados antes da primeira e da segunda prova escrita, respectivamente.\\v-a1 e a2 s\u00e3o os fatores multiplicadores das notas m\u00e9dias dos trabalhos, a serem definidos pelo docente antes do in\u00edcio de cada turma com base nas atividades espec\u00edficas a serem propostas. Os valores ser\u00e3o \u22650,1, sendo informados aos alunos no in\u00edcio do semestre. \\vEm todos os casos, os valores m\u00e1ximos para M1 e M2 ser\u00e3o \u201cdez\u201d, sendo desconsideradas pontua\u00e7\u00f5es superiores.\\vO c\u00e1lculo de NF ser\u00e1 feito conforme segue:\\vNF=(M1+2\u00d7M2)/3\\vSer\u00e3o aprovados os alunos que obtiverem NF maior ou igual 5,0.\"],\n  [\"Ser\u00e1 oferecido um programa de recupera\u00e7\u00e3o, sendo este avaliado por uma prova escrita final (PR). A m\u00e9dia de recupera\u00e7\u00e3o (MR) ser\u00e1 calculada conforme segue: MR=(NF+PR)/2Ser\u00e3o aprovados os alunos que obtiverem MR maior ou igual a 5,0.\", \"Ser\u00e1 oferecido um programa de recupera\u00e7\u00e3o, sendo este avaliado por uma prova escrita final (PR). A m\u00e9dia de recupera\u00e7\u00e3o (MR) ser\u00e1 calculada conforme segue: \\vMR=(NF+PR)/2\\vSer\u00e3o aprovados os alunos que obtiverem MR maior ou igual a 5,0.\"],\n  [\"ALTERTHUM, F.; SCHMIDELL, W.; LIMA, U. A.; MORAES. M. O. (Org.). Biotecnologia Industrial. Volume 2: Engenharia Bioqu\u00edmica. 2\u00aa Edi\u00e7\u00e3o. S\u00e3o Paulo: Blucher, 2021. p. 37-52.  ISBN 978-65-5506-019-5 (e-Book); 978-65-5506-018-8 (Impresso).DORAN P.M.; MORRISSEY, K.; CARLSON, R. P. Bioprocess Engineering Principles, 3rd edition, Academic Press, 2024. ISBN 978-0128221914SHULER, M. L.; KARGI, F.; DELISA, M. Bioprocess Engineering: Basic Concepts (3rd Edition) (Prentice Hall International Series in the Physical and Chemical Engineering Sciences) 3rd Edition. Prentice Hall; 3 edition, 2017. ISBN: 978-0137062706.\", \"ALTERTHUM, F.; SCHMIDELL, W.; LIMA, U. A.; MORAES. M. O. (Org.). Biotecnologia Industrial. Volume 2: Engenharia Bioqu\u00edmica. 2\u00aa Edi\u00e7\u00e3o. S\u00e3o Paulo: Blucher, 2021. p. 37-52.  ISBN 978-65-5506-019-5 (e-Book); 978-65-5506-018-8 (Impresso).\\vDORAN P.M.; MORRISSEY, K.; CARLSON, R. P. Bioprocess Engineering Principles, 3rd edition, Academic Press, 2024. ISBN 978-0128221914\\vSHULER, M. L.; KARGI, F.; DELISA, M. Bioprocess Engineering: Basic Concepts (3rd Edition) (Prentice Hall International Series in the Physical and Chemical Engineering Sciences) 3rd Edition. Prentice Hall; 3 edition, 2017. ISBN: 978-0137062706.\"],\n];\n\nconst body = context.document.body;\nfor (const [original, replacement] of replacements) {\n  const results = body.search(original, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Target text not found: ' + original.substring(0, 60));\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each long single-run text block with the same text split across\n# multiple runs joined by manual line breaks (Word's Shift+Enter -> OOXML <w:br/>),\n# matching the target diff's <w:t>...</w:t><w:br/><w:t>...</w:t> pattern.\n# Find/Replace's `^l` replacement code is a manual line break.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"1.The importance of oxygen transfer in biotechnological processes.2.Oxygen transfer systems.3.Oxygen transfer and microbial respiration: joint analysis of transfer and oxygen consumption, determination of kLa and QO2 in the fermentative process.4.Oxygen transfer in agitated and aerated systems: agitation of Newtonian liquids, agitation of Newtonian liquids subjected to aeration, agitation of non-Newtonian liquids, oxygen transfer.5.Scale-up variation: criteria for scaling up, comparisons between criteria for scaling up, scale-down.6.Downstream of biotechnological products: clarification, cell disruption, bioproduct separation techniques.\", \"1.The importance of oxygen transfer in biotechnological processes.^l2.Oxygen transfer systems.^l3.Oxygen transfer and microbial respiration: joint analysis of transfer and oxygen consumption, determination of kLa and QO2 in the fermentative process.^l4.Oxygen transfer in agitated and aerated systems: agitation of Newtonian liquids, agitation of Newtonian liquids subjected to aeration, agitation of non-Newtonian liquids, oxygen transfer.^l5.Scale-up variation: criteria for scaling up, comparisons between criteria for scaling up, scale-down.^l6.Downstream of biotechnological products: clarification, cell disruption, bioproduct separation techniques.\"),\n    @(\"A nota final (NF) ser\u00e1 composta pelas m\u00e9dias M1  e M2,calculadas conforme segue:M1=P1+a1\u00d7T1M2=P2+a2\u00d7T2Em que:-P1 e P2 s\u00e3o as notas da primeira e da segunda prova escrita, respectivamente (notas de zero a dez).-T1 e T2 s\u00e3o as notas m\u00e9dias dos trabalhos (notas de zero a dez) realizados antes da primeira e da segunda prova escrita, respectivamente.-a1 e a2 s\u00e3o os fatores multiplicadores das notas m\u00e9dias dos trabalhos, a serem definidos pelo docente antes do in\u00edcio de cada turma com base nas atividades espec\u00edficas a serem propostas. Os valores ser\u00e3o \u22650,1, sendo informados aos alunos no in\u00edcio do semestre. Em todos os casos, os valores m\u00e1ximos para M1 e M2 ser\u00e3o \u201cdez\u201d, sendo desconsideradas pontua\u00e7\u00f5es superiores.O c\u00e1lculo de NF ser\u00e1 feito conforme segue:NF=(M1+2\u00d7M2)/3Ser\u00e3o aprovados os alunos que obtiverem NF maior ou igual 5,0.\", \"A nota final (NF) ser\u00e1 composta pelas m\u00e9dias M1  e M2,calculadas conforme segue:^lM1=P1+a1\u00d7T1^lM2=P2+a2\u00d7T2^lEm que:^l-P1 e P2 s\u00e3o as notas da primeira e da segunda prova escrita, respectivamente (notas de zero a dez).^l-T1 e T2 s\u00e3o as notas m\u00e9dias dos trabalhos (notas de zero a dez) realizados antes da primeira e da segunda prova escrita, respectivamente.^l-a1 e a2 s\u00e3o os fatores multiplicadores das notas m\u00e9dias dos trabalhos, a serem definidos pelo docente antes do in\u00edcio de cada turma com base nas atividades espec\u00edficas a serem propostas. Os valores ser\u00e3o \u22650,1, sendo informados aos alunos no in\u00edcio do semestre. ^lEm todos os casos, os valores m\u00e1ximos para M1 e M2 ser\u00e3o \u201cdez\u201d, sendo desconsideradas pontua\u00e7\u00f5es superiores.^lO c\u00e1lculo de NF ser\u00e1 feito conforme segue:^lNF=(M1+2\u00d7M2)/3^lSer\u00e3o aprovados os alunos que obtiverem NF maior ou igual 5,0.\"),\n    @(\"Ser\u00e1 oferecido um programa de recupera\u00e7\u00e3o, sendo este avaliado por uma prova escrita final (PR). A m\u00e9dia de recupera\u00e7\u00e3o (MR) ser\u00e1 calculada conforme segue: MR=(NF+PR)/2Ser\u00e3o aprovados os alunos que obtiverem MR maior ou igual a 5,0.\", \"Ser\u00e1 oferecido um programa de recupera\u00e7\u00e3o, sendo este avaliado por uma prova escrita final (PR). A m\u00e9dia de recupera\u00e7\u00e3o (MR) ser\u00e1 calculada conforme segue: ^lMR=(NF+PR)/2^lSer\u00e3o aprovados os alunos que obtiverem MR maior ou igual a 5,0.\"),\n    @(\"ALTERTHUM, F.; SCHMIDELL, W.; LIMA, U. A.; MORAES. M. O. (Org.). Biotecnologia Industrial. Volume 2: Engenharia Bioqu\u00edmica. 2\u00aa Edi\u00e7\u00e3o. S\u00e3o Paulo: Blucher, 2021. p. 37-52.  ISBN 978-65-5506-019-5 (e-Book); 978-65-5506-018-8 (Impresso).DORAN P.M.; MORRISSEY, K.; CARLSON, R. P. Bioprocess Engineering Principles, 3rd edition, Academic Press, 2024. ISBN 978-0128221914SHULER, M. L.; KARGI, F.; DELISA, M. Bioprocess Engineering: Basic Concepts (3rd Edition) (Prentice Hall International Series in the Physical and Chemical Engineering Sciences) 3rd Edition. Prentice Hall; 3 edition, 2017. ISBN: 978-0137062706.\", \"ALTERTHUM, F.; SCHMIDELL, W.; LIMA, U. A.; MORAES. M. O. (Org.). Biotecnologia Industrial. Volume 2: Engenharia Bioqu\u00edmica. 2\u00aa Edi\u00e7\u00e3o. S\u00e3o Paulo: Blucher, 2021. p. 37-52.  ISBN 978-65-5506-019-5 (e-Book); 978-65-5506-018-8 (Impresso).^lDORAN P.M.; MORRISSEY, K.; CARLSON, R. P. Bioprocess Engineering Principles, 3rd edition, Academic Press, 2024. ISBN 978-0128221914^lSHULER, M. L.; KARGI, F.; DELISA, M. Bioprocess Engineering: Basic Concepts (3rd Edition) (Prentice Hall International Series in the Physical and Chemical Engineering Sciences) 3rd Edition. Prentice Hall; 3 edition, 2017. ISBN: 978-0137062706.\"),\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $result = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $result) {\n        throw \"Target text not found: \" + $findText.Substring(0, 60)\n    }\n}\n"}
